$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Insert two new rows into the low-stock items table, keeping the
#    list sorted alphabetically by item name (A-Z, Arabic names last):
#      - "DOLIPRANE 1 GM 15 TABS."            -> before "HI-POTENCY FORMULA 30 TABS." (row 10)
#      - "PANADOL COLD & FLU DAY 24 F.C. TABS" -> before "PHENADONE SYRUP 100 ML"       (row 13,
#        which is row 14 once the first new row has been inserted)
# ------------------------------------------------------------------
$ws.Rows.Item(10).Insert()
$ws.Rows.Item(14).Insert()

# Row heights for the newly inserted (currently blank) rows
$ws.Rows.Item(10).RowHeight = 24.75
$ws.Rows.Item(14).RowHeight = 25.5

# Re-create the merges that Insert() does not automatically add for the
# brand new rows (rows that were merely shifted down keep their merges).
$ws.Range("A10:B10").Merge()
$ws.Range("C10:G10").Merge()
$ws.Range("H10:K10").Merge()
$ws.Range("L10:M10").Merge()
$ws.Range("N10:O10").Merge()

$ws.Range("A14:B14").Merge()
$ws.Range("C14:G14").Merge()
$ws.Range("H14:K14").Merge()
$ws.Range("L14:M14").Merge()
$ws.Range("N14:O14").Merge()

# ------------------------------------------------------------------
# 2) Fill in the data for the two new rows
# ------------------------------------------------------------------
$ws.Cells.Item(10,1).Value = 4
$ws.Cells.Item(10,3).Value = "DOLIPRANE 1 GM 15 TABS."
$ws.Cells.Item(10,8).Value = "7:1"
$ws.Cells.Item(10,12).Value = "1"
$ws.Cells.Item(10,14).Value = "48.00"
$ws.Cells.Item(10,16).Value = "48.0000"
$ws.Cells.Item(10,17).Value = "1:0"

$ws.Cells.Item(14,1).Value = 8
$ws.Cells.Item(14,3).Value = "PANADOL COLD & FLU DAY 24 F.C. TABS"
$ws.Cells.Item(14,8).Value = "1:0"
$ws.Cells.Item(14,12).Value = "1"
$ws.Cells.Item(14,14).Value = "76.00"
$ws.Cells.Item(14,16).Value = "76.0000"
$ws.Cells.Item(14,17).Value = "1:0"

# ------------------------------------------------------------------
# 3) Renumber the "م" (row number) column for every item row, 1..13
# ------------------------------------------------------------------
for ($i = 0; $i -lt 13; $i++) {
  $ws.Cells.Item(7 + $i, 1).Value = $i + 1
}

# ------------------------------------------------------------------
# 4) Update the total (row 20, column P) to reflect the two new prices
# ------------------------------------------------------------------
$ws.Cells.Item(20,16).Value = 524.98

# ------------------------------------------------------------------
# 5) Update the printed timestamp footer
# ------------------------------------------------------------------
$ws.Cells.Item(21,1).Value = "Thursday, 9 October, 2025 12:35 PM"
